$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.2363035678863525
$ws.Range("E2").Value = 94.70037507160487
$ws.Range("F2").Value = 0.003912046872870529
$ws.Range("G2").Value = 0.003142746413903732
$ws.Range("H2").Value = 0.002916296946810515
$ws.Range("I2").Value = 0.00259675110709333
$ws.Range("J2").Value = 0.002497111106785256
$ws.Range("K2").Value = 0.002329977600918434
$ws.Range("L2").Value = 0.002327124338143864
$ws.Range("M2").Value = 0.00224909371341043
$ws.Range("N2").Value = 0.002120173656072495
$ws.Range("O2").Value = 0.002088404292888256
$ws.Range("P2").Value = 0.002088404292888256
$ws.Range("Q2").Value = 0.001985790224173586
$ws.Range("R2").Value = 0.001985790224173586
$ws.Range("S2").Value = 0.001942593634824211
$ws.Range("T2").Value = 0.001941730499239112
$ws.Range("U2").Value = 0.001909091485422056
$ws.Range("V2").Value = 0.001900358831564524
$ws.Range("W2").Value = 0.001848966505993615
$ws.Range("X2").Value = 0.001848966505993615
$ws.Range("Y2").Value = 0.001846011209972804

$ws.Range("C3").Value = 0.2187149524688721
$ws.Range("E3").Value = 95.20840860236967
$ws.Range("F3").Value = 0.003824426363668594
$ws.Range("G3").Value = 0.003033984057157837
$ws.Range("H3").Value = 0.002739035566543026
$ws.Range("I3").Value = 0.002691487128578803
$ws.Range("J3").Value = 0.002534468417665404
$ws.Range("K3").Value = 0.002414690657684051
$ws.Range("L3").Value = 0.002227772043741701
$ws.Range("M3").Value = 0.002227772043741701
$ws.Range("N3").Value = 0.002227772043741701
$ws.Range("O3").Value = 0.002167260317866401
$ws.Range("P3").Value = 0.00210600103979604
$ws.Range("Q3").Value = 0.002002549846542588
$ws.Range("R3").Value = 0.002002549846542588
$ws.Range("S3").Value = 0.001980474195224921
$ws.Range("T3").Value = 0.001936667856547072
$ws.Range("U3").Value = 0.001936667856547072
$ws.Range("V3").Value = 0.001927319104364032
$ws.Range("W3").Value = 0.001899896227532381
$ws.Range("X3").Value = 0.001895732526166729
$ws.Range("Y3").Value = 0.001855914397707011

$ws.Range("C4").Value = 0.2548689842224121
$ws.Range("E4").Value = 96.31415283774004
$ws.Range("F4").Value = 0.004088656908542934
$ws.Range("G4").Value = 0.003188519095289491
$ws.Range("H4").Value = 0.002826993516747557
$ws.Range("I4").Value = 0.00253370270048579
$ws.Range("J4").Value = 0.002510071378114267
$ws.Range("K4").Value = 0.002436008553278002
$ws.Range("L4").Value = 0.002307062173669521
$ws.Range("M4").Value = 0.002304529120360171
$ws.Range("N4").Value = 0.00213263695706378
$ws.Range("O4").Value = 0.00213263695706378
$ws.Range("P4").Value = 0.002044769632843683
$ws.Range("Q4").Value = 0.00201501180122241
$ws.Range("R4").Value = 0.002006823105361919
$ws.Range("S4").Value = 0.002006823105361919
$ws.Range("T4").Value = 0.001981057142960633
$ws.Range("U4").Value = 0.001961654236049007
$ws.Range("V4").Value = 0.001953096258322533
$ws.Range("W4").Value = 0.001927873118816848
$ws.Range("X4").Value = 0.001887710368261223
$ws.Range("Y4").Value = 0.001877468866232749

$ws.Range("C5").Value = 0.2595765590667725
$ws.Range("E5").Value = 94.15481710219319
$ws.Range("F5").Value = 0.004088656908542934
$ws.Range("G5").Value = 0.003167128781683876
$ws.Range("H5").Value = 0.002650386766786983
$ws.Range("I5").Value = 0.002650386766786983
$ws.Range("J5").Value = 0.002569572591692441
$ws.Range("K5").Value = 0.002471967140042116
$ws.Range("L5").Value = 0.002327637697703872
$ws.Range("M5").Value = 0.00221914969507124
$ws.Range("N5").Value = 0.00221914969507124
$ws.Range("O5").Value = 0.00214986476877929
$ws.Range("P5").Value = 0.002089829950233025
$ws.Range("Q5").Value = 0.002045623023960504
$ws.Range("R5").Value = 0.001959823391023686
$ws.Range("S5").Value = 0.001946659754929264
$ws.Range("T5").Value = 0.001925868125606896
$ws.Range("U5").Value = 0.001895647131779101
$ws.Range("V5").Value = 0.001881005564019717
$ws.Range("W5").Value = 0.001862541596680446
$ws.Range("X5").Value = 0.001862541596680446
$ws.Range("Y5").Value = 0.001835376551699672

$ws.Range("C6").Value = 0.4395735263824463
$ws.Range("E6").Value = 97.68402304501797
$ws.Range("F6").Value = 0.003827928855980839
$ws.Range("G6").Value = 0.003355105265041469
$ws.Range("H6").Value = 0.002956740630442933
$ws.Range("I6").Value = 0.002812607691111194
$ws.Range("J6").Value = 0.002633427752053182
$ws.Range("K6").Value = 0.002499946436668249
$ws.Range("L6").Value = 0.002212691053380368
$ws.Range("M6").Value = 0.002212691053380368
$ws.Range("N6").Value = 0.002144795249100991
$ws.Range("O6").Value = 0.002144795249100991
$ws.Range("P6").Value = 0.002103065095640086
$ws.Range("Q6").Value = 0.002103065095640086
$ws.Range("R6").Value = 0.002089576384692727
$ws.Range("S6").Value = 0.002027548278817695
$ws.Range("T6").Value = 0.002012183148704707
$ws.Range("U6").Value = 0.001985750255579169
$ws.Range("V6").Value = 0.00195852799214739
$ws.Range("W6").Value = 0.00195852799214739
$ws.Range("X6").Value = 0.001931394466499053
$ws.Range("Y6").Value = 0.001904171989181636

$ws.Range("C7").Value = 0.4772224426269531
$ws.Range("E7").Value = 97.425789983492
$ws.Range("F7").Value = 0.004088656908542934
$ws.Range("G7").Value = 0.003008644284585026
$ws.Range("H7").Value = 0.002796784694953538
$ws.Range("I7").Value = 0.002625403076270596
$ws.Range("J7").Value = 0.002425804352118234
$ws.Range("K7").Value = 0.002395997927293212
$ws.Range("L7").Value = 0.002332061780499098
$ws.Range("M7").Value = 0.002261076352156292
$ws.Range("N7").Value = 0.002257144189844238
$ws.Range("O7").Value = 0.002130817816102781
$ws.Range("P7").Value = 0.002119179601236167
$ws.Range("Q7").Value = 0.00206094585937911
$ws.Range("R7").Value = 0.002035854849792728
$ws.Range("S7").Value = 0.001997092587081678
$ws.Range("T7").Value = 0.001997092587081678
$ws.Range("U7").Value = 0.001983462694208097
$ws.Range("V7").Value = 0.00193242293553109
$ws.Range("W7").Value = 0.00193242293553109
$ws.Range("X7").Value = 0.001907114239150468
$ws.Range("Y7").Value = 0.001899138206305886

$ws.Range("C8").Value = 0.2981743812561035
$ws.Range("E8").Value = 89.91431406558331
$ws.Range("F8").Value = 0.003937131508040831
$ws.Range("G8").Value = 0.002849996082264698
$ws.Range("H8").Value = 0.002395031923795181
$ws.Range("I8").Value = 0.002395031923795181
$ws.Range("J8").Value = 0.002395031923795181
$ws.Range("K8").Value = 0.002190244324082543
$ws.Range("L8").Value = 0.002190244324082543
$ws.Range("M8").Value = 0.002138597953042683
$ws.Range("N8").Value = 0.002077583613136261
$ws.Range("O8").Value = 0.00199726119420721
$ws.Range("P8").Value = 0.001876410596006913
$ws.Range("Q8").Value = 0.001876410596006913
$ws.Range("R8").Value = 0.001876410596006913
$ws.Range("S8").Value = 0.001876410596006913
$ws.Range("T8").Value = 0.001866718431505553
$ws.Range("U8").Value = 0.001858388205566892
$ws.Range("V8").Value = 0.001797244979438526
$ws.Range("W8").Value = 0.001787082287252228
$ws.Range("X8").Value = 0.0017788459513857
$ws.Range("Y8").Value = 0.001752715673793046

$ws.Range("C9").Value = 0.3364613056182861
$ws.Range("E9").Value = 94.26018845775252
$ws.Range("F9").Value = 0.003961599652601022
$ws.Range("G9").Value = 0.003147208202534984
$ws.Range("H9").Value = 0.00276136930527132
$ws.Range("I9").Value = 0.002709413185633462
$ws.Range("J9").Value = 0.002582009139744958
$ws.Range("K9").Value = 0.002460449475662389
$ws.Range("L9").Value = 0.002376377655886418
$ws.Range("M9").Value = 0.002376377655886418
$ws.Range("N9").Value = 0.002362633180331415
$ws.Range("O9").Value = 0.002125105205214231
$ws.Range("P9").Value = 0.002094930800245903
$ws.Range("Q9").Value = 0.002052766269387834
$ws.Range("R9").Value = 0.002023483336081333
$ws.Range("S9").Value = 0.002023346916692382
$ws.Range("T9").Value = 0.001950313927089112
$ws.Range("U9").Value = 0.001927779234327523
$ws.Range("V9").Value = 0.001892164536685408
$ws.Range("W9").Value = 0.001888892963177047
$ws.Range("X9").Value = 0.00186780085777451
$ws.Range("Y9").Value = 0.001837430574225195

$ws.Range("C10").Value = 0.2732524871826172
$ws.Range("E10").Value = 98.45973366294675
$ws.Range("F10").Value = 0.004088656908542934
$ws.Range("G10").Value = 0.003420132961171325
$ws.Range("H10").Value = 0.002737622988815733
$ws.Range("I10").Value = 0.002737622988815733
$ws.Range("J10").Value = 0.002517870783216685
$ws.Range("K10").Value = 0.002517870783216685
$ws.Range("L10").Value = 0.002480363029723047
$ws.Range("M10").Value = 0.002327147550051953
$ws.Range("N10").Value = 0.002227114202696168
$ws.Range("O10").Value = 0.002227114202696168
$ws.Range("P10").Value = 0.002194403876291667
$ws.Range("Q10").Value = 0.00215923938622246
$ws.Range("R10").Value = 0.002126004791515649
$ws.Range("S10").Value = 0.002002562540642137
$ws.Range("T10").Value = 0.002002562540642137
$ws.Range("U10").Value = 0.002002562540642137
$ws.Range("V10").Value = 0.001950618879962856
$ws.Range("W10").Value = 0.001950618879962856
$ws.Range("X10").Value = 0.00192612236489938
$ws.Range("Y10").Value = 0.001919293053858611

$ws.Range("C11").Value = 0.2324929237365723
$ws.Range("E11").Value = 89.19542065313908
$ws.Range("F11").Value = 0.004088656908542934
$ws.Range("G11").Value = 0.003099956960999101
$ws.Range("H11").Value = 0.002996813960896625
$ws.Range("I11").Value = 0.002574166444560778
$ws.Range("J11").Value = 0.002393113738865556
$ws.Range("K11").Value = 0.002350434077441952
$ws.Range("L11").Value = 0.002265976673724578
$ws.Range("M11").Value = 0.002015734896956868
$ws.Range("N11").Value = 0.00197154275520148
$ws.Range("O11").Value = 0.00189728674122431
$ws.Range("P11").Value = 0.001892257601084317
$ws.Range("Q11").Value = 0.001856232661016248
$ws.Range("R11").Value = 0.001855731421051743
$ws.Range("S11").Value = 0.001836981112199098
$ws.Range("T11").Value = 0.001836981112199098
$ws.Range("U11").Value = 0.001799957763047011
$ws.Range("V11").Value = 0.001770140955077348
$ws.Range("W11").Value = 0.001764576782231431
$ws.Range("X11").Value = 0.001747065291152222
$ws.Range("Y11").Value = 0.001738702156981268

